$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.248339533805847
$ws.Range("B1").Value = 2.606176137924194
$ws.Range("C1").Value = 5.06770658493042
$ws.Range("D1").Value = 2.011683464050293
$ws.Range("E1").Value = 1.163261294364929
